$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.530.49'
$ws.Range('E2').Value = '  +1.03%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.378.18'
$ws.Range('E3').Value = '  +3.49%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '309.93'
$ws.Range('E5').Value = '  -0.11%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '105.21'
$ws.Range('E6').Value = '  +4.15%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.514'
$ws.Range('E7').Value = '  -4.20%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.519'
$ws.Range('E9').Value = '  -1.00%  '
$ws.Range('E10').Value = '  +0.76%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '53.39'
$ws.Range('E11').Value = '  +2.43%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0812'
$ws.Range('E12').Value = '  -1.43%  '
$ws.Range('E13').Value = '  -0.46%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.01'
$ws.Range('E14').Value = '  -1.41%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.745.69'
$ws.Range('E15').Value = '  +3.43%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.68'
$ws.Range('E16').Value = '  +4.49%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.377.65'
$ws.Range('E17').Value = '  +3.47%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.813'
$ws.Range('E18').Value = '  +0.74%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '43.489.23'
$ws.Range('E19').Value = '  +1.10%  '
$ws.Range('E20').Value = '  -4.15%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.30'
$ws.Range('E21').Value = '  +3.77%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0₃0922'
$ws.Range('E22').Value = '  -0.34%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '68.45'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '241.78'
$ws.Range('E24').Value = '  +0.57%  '
$ws.Range('E25').Value = '  +2.55%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.63'
$ws.Range('E26').Value = '  +0.09%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.999'
$ws.Range('E27').Value = '  -0.03%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '25.76'
$ws.Range('E28').Value = '  +5.26%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.86'
$ws.Range('E29').Value = '  -2.93%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.29'
$ws.Range('E30').Value = '  +3.57%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '37.01'
$ws.Range('E31').Value = '  -3.90%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '9.55'
$ws.Range('E32').Value = '  -0.98%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '162.10'
$ws.Range('E33').Value = '  -3.28%  '
$ws.Range('E34').Value = '  -0.88%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '18.44'
$ws.Range('E35').Value = '  +3.78%  '
$ws.Range('E36').Value = '  +0.00%  '
$ws.Range('E37').Value = '  +6.16%  '
$ws.Range('B38').Value = 'LidoDAOToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.10'
$ws.Range('E38').Value = '  -0.70%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.69'
$ws.Range('E39').Value = '  +11.10%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0743'
$ws.Range('E40').Value = '  +0.31%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.95'
$ws.Range('E42').Value = '  -1.46%  '
$ws.Range('E43').Value = '  -1.77%  '
$ws.Range('B44').Value = 'ApeXProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.68'
$ws.Range('E44').Value = '  +16.94%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '19.86'
$ws.Range('E45').Value = '  +4.36%  '
$ws.Range('B46').Value = 'Maker'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.027.42'
$ws.Range('E46').Value = '  +2.72%  '
$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0290'
$ws.Range('E47').Value = '  +0.31%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.16'
$ws.Range('E48').Value = '  +4.00%  '
$ws.Range('B49').Value = 'FraxShare'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '10.62'
$ws.Range('E49').Value = '  +7.68%  '
$ws.Range('B50').Value = 'MultiversX'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '58.12'
$ws.Range('E50').Value = '  +4.01%  '
$ws.Range('B51').Value = 'HuobiToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.97'
$ws.Range('E51').Value = '  +0.78%  '
